{"js": "// Circle Language Spec Plan: Set font to Calibri for non-heading text.\n\n// 1. The body (non-heading) text uses the \"Normal\" paragraph style -\n//    switch its font from Tahoma to Calibri and give it an explicit\n//    11pt (22 half-point) size.\nconst styles = context.document.getStyles();\nconst normalStyle = styles.getByNameOrNullObject(\"Normal\");\nnormalStyle.load(\"isNullObject\");\nawait context.sync();\n\nif (!normalStyle.isNullObject) {\n  normalStyle.font.name = \"Calibri\";\n  normalStyle.font.size = 11;\n  await context.sync();\n}\n\n// 2. Word stamps a \"_GoBack\" bookmark at the location of the most\n//    recent edit. It previously sat at the very top of the document;\n//    remove it there and re-create it where the last edit was made -\n//    right after the word \"initial\" in the \"different document...\"\n//    remark paragraph (splitting that run in two around the mark).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst body = context.document.body;\nconst matches = body.search(\"different document. But the initial\", { matchCase: true });\nmatches.load(\"text\");\nawait context.sync();\n\nif (matches.items.length > 0) {\n  const target = matches.items[0];\n  const splitPoint = target.getRange(\"End\");\n  splitPoint.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Circle Language Spec Plan: Set font to Calibri for non-heading text.\n$d = $word.ActiveDocument\n\n# 1. The body (non-heading) text style is \"Normal\" - switch its font from\n#    Tahoma to Calibri and give it an explicit 11pt (22 half-point) size.\n$normal = $d.Styles(\"Normal\")\n$normal.Font.Name = \"Calibri\"\n$normal.Font.Size = 11\n\n# 2. Word stamps a \"_GoBack\" bookmark at the location of the most recent\n#    edit. Previously it sat at the very top of the document; move it to\n#    mark the spot where the last edit was made, right after the word\n#    \"initial\" in the \"different document...\" remark paragraph.\n#    (\"_GoBack\" is a hidden bookmark, so Bookmarks.Count/.Exists don't see\n#    it - index it directly and swallow the error if it is already gone.)\ntry {\n    $d.Bookmarks(\"_GoBack\").Delete()\n} catch {\n}\n\n$searchRange = $d.Content\n$find = $searchRange.Find\n$find.Text = \"different document. But the initial\"\n$find.Execute() | Out-Null\n\n$splitPoint = $searchRange.End\n$bmRange = $d.Range($splitPoint, $splitPoint)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange)\n"}
